$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '95.704.94'
$ws.Range('E2').Value = '  +3.16%  '

$ws.Range('D3').Value = '3.598.25'
$ws.Range('E3').Value = '  +5.64%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').Value = '''239.18'
$ws.Range('E5').Value = '  +4.21%  '

$ws.Range('D6').Value = '''655.05'
$ws.Range('E6').Value = '  +6.04%  '

$ws.Range('E7').Value = '  +7.52%  '

$ws.Range('D8').Value = '''0.409'
$ws.Range('E8').Value = '  +4.97%  '

$ws.Range('E9').Value = '  -0.09%  '

$ws.Range('E10').Value = '  +5.32%  '

$ws.Range('D11').Value = '3.595.87'
$ws.Range('E11').Value = '  +5.55%  '

$ws.Range('D12').Value = '''43.29'
$ws.Range('E12').Value = '  +1.28%  '

$ws.Range('E13').Value = '  +2.01%  '

$ws.Range('D14').Value = '''6.34'
$ws.Range('E14').Value = '  +1.88%  '

$ws.Range('D15').Value = '4.271.18'
$ws.Range('E15').Value = '  +5.85%  '

$ws.Range('D16').Value = '95.489.88'
$ws.Range('E16').Value = '  +3.09%  '

$ws.Range('D17').Value = '''0.0000257'
$ws.Range('E17').Value = '  +4.76%  '

$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.605.78'
$ws.Range('E18').Value = '  +5.88%  '

$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').Value = '''7.95'
$ws.Range('E19').Value = '  -1.30%  '

$ws.Range('D20').Value = '''12.55'
$ws.Range('E20').Value = '  +9.30%  '

$ws.Range('D21').Value = '''18.10'
$ws.Range('E21').Value = '  +1.45%  '

$ws.Range('D22').Value = '''3.55'
$ws.Range('E22').Value = '  +7.39%  '

$ws.Range('E23').Value = '  +12.38%  '

$ws.Range('D24').Value = '''512.16'
$ws.Range('E24').Value = '  +3.51%  '

$ws.Range('D25').Value = '''0.0000197'
$ws.Range('E25').Value = '  +7.40%  '

$ws.Range('D26').Value = '''6.65'
$ws.Range('E26').Value = '  +2.11%  '

$ws.Range('D27').Value = '''97.05'
$ws.Range('E27').Value = '  +2.86%  '

$ws.Range('D28').Value = '''12.83'
$ws.Range('E28').Value = '  +8.13%  '

$ws.Range('D29').Value = '3.800.35'
$ws.Range('E29').Value = '  +6.06%  '

$ws.Range('D30').Value = '''3.24'
$ws.Range('E30').Value = '  +19.43%  '

$ws.Range('D31').Value = '''11.37'
$ws.Range('E31').Value = '  +1.16%  '

$ws.Range('D32').Value = '''0.998'
$ws.Range('E32').Value = '  -0.24%  '

$ws.Range('E33').Value = '  +4.01%  '

$ws.Range('D34').Value = '''1.01'
$ws.Range('E34').Value = '  +0.74%  '

$ws.Range('D35').Value = '''0.177'
$ws.Range('E35').Value = '  +3.06%  '

$ws.Range('D36').Value = '''31.92'
$ws.Range('E36').Value = '  +7.45%  '

$ws.Range('D37').Value = '''0.562'
$ws.Range('E37').Value = '  +4.83%  '

$ws.Range('D38').Value = '''8.27'
$ws.Range('E38').Value = '  +11.54%  '

$ws.Range('D39').Value = '''566.34'
$ws.Range('E39').Value = '  +3.11%  '

$ws.Range('D40').Value = '''1.50'
$ws.Range('E40').Value = '  +8.17%  '

$ws.Range('D41').Value = '''0.151'
$ws.Range('E41').Value = '  +1.68%  '

$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').Value = '''1.00'
$ws.Range('E42').Value = '  -0.04%  '

$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = '''0.928'
$ws.Range('E43').Value = '  +2.19%  '

$ws.Range('D44').Value = '''1.73'
$ws.Range('E44').Value = '  +1.41%  '

$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').Value = '''5.74'
$ws.Range('E45').Value = '  +5.50%  '

$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').Value = '''23.77'
$ws.Range('E46').Value = '  +0.51%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '''34.10'
$ws.Range('E47').Value = '  +35.90%  '

$ws.Range('D48').Value = '''0.0419'
$ws.Range('E48').Value = '  +3.94%  '

$ws.Range('E49').Value = '  +7.96%  '

$ws.Range('D50').Value = '''54.37'
$ws.Range('E50').Value = '  +2.10%  '

$ws.Range('D51').Value = '''3.46'
$ws.Range('E51').Value = '  -5.53%  '
